$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad / Changed date) for rows 2-13 from 45174 to 45175
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
